# Fix cost bug: update calibrated values in columns J:AS for specific rows
# on the single worksheet "strategy_id-0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    96  = 9288329.044
    97  = 1856357.96
    98  = 744216.8137000001
    99  = 54560.23073
    100 = 1680518.694
    101 = 15173145.3
    103 = 1423204.72
    104 = 752788.9917
    111 = 52831.46011
    112 = 214809.8038
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $range = $ws.Range("J$row`:AS$row")
    $range.Value = $value
}
